$wb = $excel.ActiveWorkbook
$ws1 = $wb.Sheets.Item("Parameter Values")
$ws2 = $wb.Sheets.Item("Results")

# --- Results sheet: add new explanation row for "G2 with m2v" ---
# Create the new label cell first (so its shared string is appended before
# the two new numeric-text values below), matching header style of A1/G1/M1.
$ws2.Range("A1").Copy()
$ws2.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("A18").Value = "for G2 with m2v"

# --- Parameter Values sheet: new metapath2vec / GNN learning-rate results for G2 ---
# N9 (lr_m2v for G2) and N17 (lr for G2, GNN) are stored as text (shared
# string) cells, not numbers, so force text entry with a leading apostrophe
# then restore the original (non quote-prefixed) number formatting by
# pasting the format from an equivalent neighboring text cell.
$ws1.Range("N9").Value = "'0.010631597403622543"
$ws1.Range("N8").Copy()
$ws1.Range("N9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("N17").Value = "'0.04801629135272712"
$ws1.Range("M17").Copy()
$ws1.Range("N17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Plain numeric hyperparameter updates for the G2 column
$ws1.Range("N5").Value = 2
$ws1.Range("N7").Value = 7
$ws1.Range("N14").Value = 256
$ws1.Range("N15").Value = 128
$ws1.Range("N16").Value = 150
$ws1.Range("N19").Value = 0

# --- Selection / active sheet state ---
# Results sheet keeps selection at B18 (new row) but is no longer the active tab.
$ws2.Range("B18").Select()

# Parameter Values becomes the active/selected tab with the cursor at N20.
$ws1.Activate()
$ws1.Range("N20").Select()
